$wb = $excel.ActiveWorkbook

# --- Update selection (active cell) on sheet "A" ---
$wsA = $wb.Worksheets.Item("A")
$wsA.Activate() | Out-Null
$wsA.Range("D8").Select() | Out-Null

# --- Update selection (active cell) on sheet "C" ---
$wsC = $wb.Worksheets.Item("C")
$wsC.Activate() | Out-Null
$wsC.Range("D4").Select() | Out-Null

# --- Insert new worksheet "After Loop" right before the "DATA" sheet ---
$dataSheet = $wb.Worksheets.Item("DATA")
$newSheet = $wb.Worksheets.Add($dataSheet)
$newSheet.Name = "After Loop"

# Populate the new "After Loop" sheet, mirroring the loop-body sheets (A/B/C)
# but with two extra calculated columns (Calculated2 / Calculated3)
$newSheet.Range("A2").Value = "{{#each items}}{{value}}"
$newSheet.Range("B2").Value = "{{/each}}"
$newSheet.Range("A2:B2").Font.Color = 0

$newSheet.Range("B4").Value = "{{#each items}}"

$newSheet.Range("C5").Value = "Value"
$newSheet.Range("D5").Value = "Calculated"
$newSheet.Range("C5:D5").Font.Bold = $true

$newSheet.Range("E5").Value = "Calculated2"
$newSheet.Range("E5").Font.Bold = $true

$newSheet.Range("F5").Value = "Calculated3"
$newSheet.Range("F5").Font.Bold = $true
$newSheet.Range("F5").Font.Color = 0

$newSheet.Range("C6").Value = "{{value}}"
$newSheet.Range("D6").Formula = "=C6+DATA!A1"
$newSheet.Range("E6").Formula = "=C6+DATA!B4"
$newSheet.Range("F6").Formula = "=C6+DATA!C7"

$newSheet.Range("B8").Value = "{{/each}}"

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

$newSheet.Activate() | Out-Null
$newSheet.Range("C7").Select() | Out-Null

# --- Update selection (active cell) on sheet "DATA" ---
$wsData = $wb.Worksheets.Item("DATA")
$wsData.Activate() | Out-Null
$wsData.Range("C17").Select() | Out-Null

# Leave "After Loop" as the active tab (matches activeTab index in the workbook)
$newSheet.Activate() | Out-Null
